$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 682.1579
$ws.Range("I6").Value = 236.75
$ws.Range("J6").Value = 1445.7142
$ws.Range("K6").Value = 710.25
$ws.Range("L6").Value = 4337.142599999999
$ws.Range("M6").Value = -598.25
$ws.Range("N6").Value = -4561.142599999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1376.56
$ws.Range("J112").Value = 1399.5464
$ws.Range("L112").Value = 4198.6392
$ws.Range("N112").Value = -6414.6392

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3033744.8
$ws.Range("I137").Value = 4549242
$ws.Range("K137").Value = 13647726
$ws.Range("M137").Value = -13645176

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2283.1125
$ws.Range("J138").Value = 5312.8076
$ws.Range("L138").Value = 15938.4228
$ws.Range("N138").Value = -26218.4228

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 14000
$ws.Range("I31").Value = 14000
$ws.Range("K31").Value = 14000
$ws.Range("M31").Value = -13706

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1483.5536
$ws.Range("I61").Value = 731.093
$ws.Range("J61").Value = 3972.4614
$ws.Range("K61").Value = 731.093
$ws.Range("L61").Value = 3972.4614
$ws.Range("M61").Value = -519.093
$ws.Range("N61").Value = -4396.4614

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 933.90625
$ws.Range("I74").Value = 750.1852
$ws.Range("J74").Value = 1926
$ws.Range("K74").Value = 750.1852
$ws.Range("L74").Value = 1926
$ws.Range("M74").Value = 123.8148
$ws.Range("N74").Value = -3674

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 933.90625
$ws.Range("I77").Value = 750.1852
$ws.Range("J77").Value = 1926
$ws.Range("K77").Value = 3750.926
$ws.Range("L77").Value = 9630
$ws.Range("M77").Value = 617.0740000000001
$ws.Range("N77").Value = -18366

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1483.5536
$ws.Range("I136").Value = 731.093
$ws.Range("J136").Value = 3972.4614
$ws.Range("K136").Value = 2193.279
$ws.Range("L136").Value = 11917.3842
$ws.Range("M136").Value = 356.721
$ws.Range("N136").Value = -17017.3842

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1420.2903
$ws.Range("I105").Value = 1389.3125
$ws.Range("J105").Value = 1453.3334
$ws.Range("K105").Value = 1389.3125
$ws.Range("L105").Value = 1453.3334
$ws.Range("M105").Value = 357.6875
$ws.Range("N105").Value = -4947.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1328.6
$ws.Range("I134").Value = 904.925
$ws.Range("J134").Value = 4718
$ws.Range("K134").Value = 2714.775
$ws.Range("L134").Value = 14154
$ws.Range("M134").Value = -179.7749999999996
$ws.Range("N134").Value = -19224

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2002747.5
$ws.Range("I31").Value = 3334952.2
$ws.Range("J31").Value = 4440.65
$ws.Range("K31").Value = 3334952.2
$ws.Range("L31").Value = 4440.65
$ws.Range("M31").Value = -3334657.2
$ws.Range("N31").Value = -5030.65

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 2002747.5
$ws.Range("I34").Value = 3334952.2
$ws.Range("J34").Value = 4440.65
$ws.Range("K34").Value = 3334952.2
$ws.Range("L34").Value = 4440.65
$ws.Range("M34").Value = -3334750.2
$ws.Range("N34").Value = -4844.65

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 8622976
$ws.Range("I58").Value = 1369.9778
$ws.Range("J58").Value = 38466996
$ws.Range("K58").Value = 1369.9778
$ws.Range("L58").Value = 38466996
$ws.Range("M58").Value = -1166.9778
$ws.Range("N58").Value = -38467402

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2264.5454
$ws.Range("I105").Value = 2008.5714
$ws.Range("J105").Value = 2712.5
$ws.Range("K105").Value = 2008.5714
$ws.Range("L105").Value = 2712.5
$ws.Range("M105").Value = -261.5714
$ws.Range("N105").Value = -6206.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1300.6
$ws.Range("I132").Value = 1033.6522
$ws.Range("J132").Value = 1946.8948
$ws.Range("K132").Value = 3100.9566
$ws.Range("L132").Value = 5840.6844
$ws.Range("M132").Value = -570.9566
$ws.Range("N132").Value = -10900.6844

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1291.5962
$ws.Range("I134").Value = 799.2449
$ws.Range("K134").Value = 2397.7347
$ws.Range("M134").Value = 137.2653

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 8622976
$ws.Range("I136").Value = 1369.9778
$ws.Range("J136").Value = 38466996
$ws.Range("K136").Value = 4109.9334
$ws.Range("L136").Value = 115400988
$ws.Range("M136").Value = -1559.9334
$ws.Range("N136").Value = -115406088

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H35").Value = 2496
$ws.Range("I35").Value = 1166.3334
$ws.Range("J35").Value = 3493.25
$ws.Range("K35").Value = 3499.0002
$ws.Range("L35").Value = 10479.75
$ws.Range("M35").Value = -3211.0002
$ws.Range("N35").Value = -11055.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 27582.5
$ws.Range("I129").Value = 2688.4614
$ws.Range("J129").Value = 73814.28999999999
$ws.Range("K129").Value = 8065.3842
$ws.Range("L129").Value = 221442.87
$ws.Range("M129").Value = -3065.3842
$ws.Range("N129").Value = -231442.87

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H130").Value = 2206
$ws.Range("J130").Value = 2500
$ws.Range("L130").Value = 7500
$ws.Range("N130").Value = -17540

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1842.3438
$ws.Range("J131").Value = 1407.2916
$ws.Range("L131").Value = 4221.8748
$ws.Range("N131").Value = -14301.8748

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 2398.75
$ws.Range("I136").Value = 1531.6666
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 4594.9998
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = 505.0002000000004
$ws.Range("N136").Value = -25200

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 6948944
$ws.Range("I139").Value = 20836860
$ws.Range("J139").Value = 4986.0415
$ws.Range("K139").Value = 62510580
$ws.Range("L139").Value = 14958.1245
$ws.Range("M139").Value = -62505440
$ws.Range("N139").Value = -25238.1245

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 9807977
$ws.Range("I140").Value = 20834152
$ws.Range("J140").Value = 6933.3335
$ws.Range("K140").Value = 62502456
$ws.Range("L140").Value = 20800.0005
$ws.Range("M140").Value = -62497276
$ws.Range("N140").Value = -31160.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H108").Value = 28684
$ws.Range("J108").Value = 28684
$ws.Range("L108").Value = 28684
$ws.Range("N108").Value = -36364

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1860.9048
$ws.Range("I132").Value = 1216.909
$ws.Range("J132").Value = 4222.222
$ws.Range("K132").Value = 3650.727
$ws.Range("L132").Value = 12666.666
$ws.Range("M132").Value = -1120.727
$ws.Range("N132").Value = -17726.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1820280.5
$ws.Range("J136").Value = 7487.375
$ws.Range("L136").Value = 22462.125
$ws.Range("N136").Value = -27562.125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 402749.03
$ws.Range("I122").Value = 771462.0600000001
$ws.Range("J122").Value = 3309.9167
$ws.Range("K122").Value = 2314386.18
$ws.Range("L122").Value = 9929.750100000001
$ws.Range("M122").Value = -2311936.18
$ws.Range("N122").Value = -14829.7501

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 200412
$ws.Range("I132").Value = 288182.88
$ws.Range("J132").Value = 29746.389
$ws.Range("K132").Value = 864548.64
$ws.Range("L132").Value = 89239.167
$ws.Range("M132").Value = -862018.64
$ws.Range("N132").Value = -94299.167

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 724.40625
$ws.Range("I136").Value = 337.1087
$ws.Range("J136").Value = 1714.1666
$ws.Range("K136").Value = 1011.3261
$ws.Range("L136").Value = 5142.4998
$ws.Range("M136").Value = 1538.6739
$ws.Range("N136").Value = -10242.4998
